$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Myoc"
$ws.Cells.Item(2, 3).Value2 = "Fzd4"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = [double]"1"
$ws.Cells.Item(2, 6).Value2 = [double]"0.3333333333333333"
$ws.Cells.Item(2, 7).Value2 = [double]"0.1603853333333333"
$ws.Cells.Item(2, 8).Value2 = [double]"0.481156"
$ws.Cells.Item(2, 9).Value2 = [double]"0.01032935781992836"
$ws.Cells.Item(2, 10).Value2 = [double]"0.01042870175281933"
$ws.Cells.Item(2, 11).Value2 = [double]"3"
$ws.Cells.Item(2, 12).Value2 = [double]"1"
$ws.Cells.Item(2, 13).Value2 = [double]"19.524618"
$ws.Cells.Item(2, 14).Value2 = [double]"58.573854"
$ws.Cells.Item(2, 15).Value2 = [double]"0.4154885426712971"
$ws.Cells.Item(2, 16).Value2 = [double]"0.4539723485554654"
$ws.Cells.Item(2, 17).Value2 = [double]"3.131462366136"
$ws.Cells.Item(2, 18).Value2 = [double]"28.183161295224"
$ws.Cells.Item(2, 19).Value2 = [double]"0.004291729827332402"
$ws.Cells.Item(2, 20).Value2 = [double]"0.004734342227111893"

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Myoc"
$ws.Cells.Item(3, 3).Value2 = "Fzd4"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = [double]"1"
$ws.Cells.Item(3, 6).Value2 = [double]"0.3333333333333333"
$ws.Cells.Item(3, 7).Value2 = [double]"0.1603853333333333"
$ws.Cells.Item(3, 8).Value2 = [double]"0.481156"
$ws.Cells.Item(3, 9).Value2 = [double]"0.01032935781992836"
$ws.Cells.Item(3, 10).Value2 = [double]"0.01042870175281933"
$ws.Cells.Item(3, 11).Value2 = [double]"3"
$ws.Cells.Item(3, 12).Value2 = [double]"1"
$ws.Cells.Item(3, 13).Value2 = [double]"15.24435933333334"
$ws.Cells.Item(3, 14).Value2 = [double]"45.73307800000001"
$ws.Cells.Item(3, 15).Value2 = [double]"0.324403614112412"
$ws.Cells.Item(3, 16).Value2 = [double]"0.3544508583357054"
$ws.Cells.Item(3, 17).Value2 = [double]"2.444971653129778"
$ws.Cells.Item(3, 18).Value2 = [double]"22.004744878168"
$ws.Cells.Item(3, 19).Value2 = [double]"0.003350881008245065"
$ws.Cells.Item(3, 20).Value2 = [double]"0.003696462287613889"

# Row 4
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Myoc"
$ws.Cells.Item(4, 3).Value2 = "Fzd4"
$ws.Cells.Item(4, 4).Value2 = "MuSCs"
$ws.Cells.Item(4, 5).Value2 = [double]"1"
$ws.Cells.Item(4, 6).Value2 = [double]"0.3333333333333333"
$ws.Cells.Item(4, 7).Value2 = [double]"0.1603853333333333"
$ws.Cells.Item(4, 8).Value2 = [double]"0.481156"
$ws.Cells.Item(4, 9).Value2 = [double]"0.01032935781992836"
$ws.Cells.Item(4, 10).Value2 = [double]"0.01042870175281933"
$ws.Cells.Item(4, 11).Value2 = [double]"2"
$ws.Cells.Item(4, 12).Value2 = [double]"1"
$ws.Cells.Item(4, 13).Value2 = [double]"11.9507005"
$ws.Cells.Item(4, 14).Value2 = [double]"23.901401"
$ws.Cells.Item(4, 15).Value2 = [double]"0.2543137660693869"
$ws.Cells.Item(4, 16).Value2 = [double]"0.1852460510065796"
$ws.Cells.Item(4, 17).Value2 = [double]"1.916717083259333"
$ws.Cells.Item(4, 18).Value2 = [double]"11.500302499556"
$ws.Cells.Item(4, 19).Value2 = [double]"0.002626897888264254"
$ws.Cells.Item(4, 20).Value2 = [double]"0.001931875816835177"

# Row 5
$ws.Cells.Item(5, 1).Value2 = "ECs"
$ws.Cells.Item(5, 2).Value2 = "Myoc"
$ws.Cells.Item(5, 3).Value2 = "Fzd4"
$ws.Cells.Item(5, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value2 = [double]"1"
$ws.Cells.Item(5, 6).Value2 = [double]"0.3333333333333333"
$ws.Cells.Item(5, 7).Value2 = [double]"0.1603853333333333"
$ws.Cells.Item(5, 8).Value2 = [double]"0.481156"
$ws.Cells.Item(5, 9).Value2 = [double]"0.01032935781992836"
$ws.Cells.Item(5, 10).Value2 = [double]"0.01042870175281933"
$ws.Cells.Item(5, 11).Value2 = [double]"2"
$ws.Cells.Item(5, 12).Value2 = [double]"0.6666666666666666"
$ws.Cells.Item(5, 13).Value2 = [double]"0.272275"
$ws.Cells.Item(5, 14).Value2 = [double]"0.816825"
$ws.Cells.Item(5, 15).Value2 = [double]"0.005794077146903843"
$ws.Cells.Item(5, 16).Value2 = [double]"0.006330742102249548"
$ws.Cells.Item(5, 17).Value2 = [double]"0.04366891663333333"
$ws.Cells.Item(5, 18).Value2 = [double]"0.3930202497"
$ws.Cells.Item(5, 19).Value2 = [double]"5.984909608663941e-05"
$ws.Cells.Item(5, 20).Value2 = [double]"6.602142125837701e-05"

# Row 6
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Myoc"
$ws.Cells.Item(6, 3).Value2 = "Fzd4"
$ws.Cells.Item(6, 4).Value2 = "ECs"
$ws.Cells.Item(6, 5).Value2 = [double]"3"
$ws.Cells.Item(6, 6).Value2 = [double]"1"
$ws.Cells.Item(6, 7).Value2 = [double]"14.92301466666667"
$ws.Cells.Item(6, 8).Value2 = [double]"44.76904399999999"
$ws.Cells.Item(6, 9).Value2 = [double]"0.9610926076617912"
$ws.Cells.Item(6, 10).Value2 = [double]"0.970336039943066"
$ws.Cells.Item(6, 11).Value2 = [double]"3"
$ws.Cells.Item(6, 12).Value2 = [double]"1"
$ws.Cells.Item(6, 13).Value2 = [double]"19.524618"
$ws.Cells.Item(6, 14).Value2 = [double]"58.573854"
$ws.Cells.Item(6, 15).Value2 = [double]"0.4154885426712971"
$ws.Cells.Item(6, 16).Value2 = [double]"0.4539723485554654"
$ws.Cells.Item(6, 17).Value2 = [double]"291.366160775064"
$ws.Cells.Item(6, 18).Value2 = [double]"2622.295446975575"
$ws.Cells.Item(6, 19).Value2 = [double]"0.3993229669295544"
$ws.Cells.Item(6, 20).Value2 = [double]"0.4405057309409636"

# Row 7
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Myoc"
$ws.Cells.Item(7, 3).Value2 = "Fzd4"
$ws.Cells.Item(7, 4).Value2 = "FAPs"
$ws.Cells.Item(7, 5).Value2 = [double]"3"
$ws.Cells.Item(7, 6).Value2 = [double]"1"
$ws.Cells.Item(7, 7).Value2 = [double]"14.92301466666667"
$ws.Cells.Item(7, 8).Value2 = [double]"44.76904399999999"
$ws.Cells.Item(7, 9).Value2 = [double]"0.9610926076617912"
$ws.Cells.Item(7, 10).Value2 = [double]"0.970336039943066"
$ws.Cells.Item(7, 11).Value2 = [double]"3"
$ws.Cells.Item(7, 12).Value2 = [double]"1"
$ws.Cells.Item(7, 13).Value2 = [double]"15.24435933333334"
$ws.Cells.Item(7, 14).Value2 = [double]"45.73307800000001"
$ws.Cells.Item(7, 15).Value2 = [double]"0.324403614112412"
$ws.Cells.Item(7, 16).Value2 = [double]"0.3544508583357054"
$ws.Cells.Item(7, 17).Value2 = [double]"227.4917979152702"
$ws.Cells.Item(7, 18).Value2 = [double]"2047.426181237432"
$ws.Cells.Item(7, 19).Value2 = [double]"0.3117819154222075"
$ws.Cells.Item(7, 20).Value2 = [double]"0.3439364422318891"

# Row 8
$ws.Cells.Item(8, 1).Value2 = "FAPs"
$ws.Cells.Item(8, 2).Value2 = "Myoc"
$ws.Cells.Item(8, 3).Value2 = "Fzd4"
$ws.Cells.Item(8, 4).Value2 = "MuSCs"
$ws.Cells.Item(8, 5).Value2 = [double]"3"
$ws.Cells.Item(8, 6).Value2 = [double]"1"
$ws.Cells.Item(8, 7).Value2 = [double]"14.92301466666667"
$ws.Cells.Item(8, 8).Value2 = [double]"44.76904399999999"
$ws.Cells.Item(8, 9).Value2 = [double]"0.9610926076617912"
$ws.Cells.Item(8, 10).Value2 = [double]"0.970336039943066"
$ws.Cells.Item(8, 11).Value2 = [double]"2"
$ws.Cells.Item(8, 12).Value2 = [double]"1"
$ws.Cells.Item(8, 13).Value2 = [double]"11.9507005"
$ws.Cells.Item(8, 14).Value2 = [double]"23.901401"
$ws.Cells.Item(8, 15).Value2 = [double]"0.2543137660693869"
$ws.Cells.Item(8, 16).Value2 = [double]"0.1852460510065796"
$ws.Cells.Item(8, 17).Value2 = [double]"178.3404788384406"
$ws.Cells.Item(8, 18).Value2 = [double]"1070.042873030644"
$ws.Cells.Item(8, 19).Value2 = [double]"0.2444190805959178"
$ws.Cells.Item(8, 20).Value2 = [double]"0.1797509195488157"

# Row 9
$ws.Cells.Item(9, 1).Value2 = "FAPs"
$ws.Cells.Item(9, 2).Value2 = "Myoc"
$ws.Cells.Item(9, 3).Value2 = "Fzd4"
$ws.Cells.Item(9, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value2 = [double]"3"
$ws.Cells.Item(9, 6).Value2 = [double]"1"
$ws.Cells.Item(9, 7).Value2 = [double]"14.92301466666667"
$ws.Cells.Item(9, 8).Value2 = [double]"44.76904399999999"
$ws.Cells.Item(9, 9).Value2 = [double]"0.9610926076617912"
$ws.Cells.Item(9, 10).Value2 = [double]"0.970336039943066"
$ws.Cells.Item(9, 11).Value2 = [double]"2"
$ws.Cells.Item(9, 12).Value2 = [double]"0.6666666666666666"
$ws.Cells.Item(9, 13).Value2 = [double]"0.272275"
$ws.Cells.Item(9, 14).Value2 = [double]"0.816825"
$ws.Cells.Item(9, 15).Value2 = [double]"0.005794077146903843"
$ws.Cells.Item(9, 16).Value2 = [double]"0.006330742102249548"
$ws.Cells.Item(9, 17).Value2 = [double]"4.063163818366666"
$ws.Cells.Item(9, 18).Value2 = [double]"36.5684743653"
$ws.Cells.Item(9, 19).Value2 = [double]"0.005568644714111406"
$ws.Cells.Item(9, 20).Value2 = [double]"0.006142947221397667"

# Row 10
$ws.Cells.Item(10, 1).Value2 = "MuSCs"
$ws.Cells.Item(10, 2).Value2 = "Myoc"
$ws.Cells.Item(10, 3).Value2 = "Fzd4"
$ws.Cells.Item(10, 4).Value2 = "ECs"
$ws.Cells.Item(10, 5).Value2 = [double]"1"
$ws.Cells.Item(10, 6).Value2 = [double]"0.5"
$ws.Cells.Item(10, 7).Value2 = [double]"0.443735"
$ws.Cells.Item(10, 8).Value2 = [double]"0.88747"
$ws.Cells.Item(10, 9).Value2 = [double]"0.02857803451828042"
$ws.Cells.Item(10, 10).Value2 = [double]"0.01923525830411462"
$ws.Cells.Item(10, 11).Value2 = [double]"3"
$ws.Cells.Item(10, 12).Value2 = [double]"1"
$ws.Cells.Item(10, 13).Value2 = [double]"19.524618"
$ws.Cells.Item(10, 14).Value2 = [double]"58.573854"
$ws.Cells.Item(10, 15).Value2 = [double]"0.4154885426712971"
$ws.Cells.Item(10, 16).Value2 = [double]"0.4539723485554654"
$ws.Cells.Item(10, 17).Value2 = [double]"8.66375636823"
$ws.Cells.Item(10, 18).Value2 = [double]"51.98253820938"
$ws.Cells.Item(10, 19).Value2 = [double]"0.01187384591441036"
$ws.Cells.Item(10, 20).Value2 = [double]"0.008732275387389936"

# Row 11
$ws.Cells.Item(11, 1).Value2 = "MuSCs"
$ws.Cells.Item(11, 2).Value2 = "Myoc"
$ws.Cells.Item(11, 3).Value2 = "Fzd4"
$ws.Cells.Item(11, 4).Value2 = "FAPs"
$ws.Cells.Item(11, 5).Value2 = [double]"1"
$ws.Cells.Item(11, 6).Value2 = [double]"0.5"
$ws.Cells.Item(11, 7).Value2 = [double]"0.443735"
$ws.Cells.Item(11, 8).Value2 = [double]"0.88747"
$ws.Cells.Item(11, 9).Value2 = [double]"0.02857803451828042"
$ws.Cells.Item(11, 10).Value2 = [double]"0.01923525830411462"
$ws.Cells.Item(11, 11).Value2 = [double]"3"
$ws.Cells.Item(11, 12).Value2 = [double]"1"
$ws.Cells.Item(11, 13).Value2 = [double]"15.24435933333334"
$ws.Cells.Item(11, 14).Value2 = [double]"45.73307800000001"
$ws.Cells.Item(11, 15).Value2 = [double]"0.324403614112412"
$ws.Cells.Item(11, 16).Value2 = [double]"0.3544508583357054"
$ws.Cells.Item(11, 17).Value2 = [double]"6.764455788776668"
$ws.Cells.Item(11, 18).Value2 = [double]"40.58673473266001"
$ws.Cells.Item(11, 19).Value2 = [double]"0.00927081768195943"
$ws.Cells.Item(11, 20).Value2 = [double]"0.006817953816202434"

# Row 12
$ws.Cells.Item(12, 1).Value2 = "MuSCs"
$ws.Cells.Item(12, 2).Value2 = "Myoc"
$ws.Cells.Item(12, 3).Value2 = "Fzd4"
$ws.Cells.Item(12, 4).Value2 = "MuSCs"
$ws.Cells.Item(12, 5).Value2 = [double]"1"
$ws.Cells.Item(12, 6).Value2 = [double]"0.5"
$ws.Cells.Item(12, 7).Value2 = [double]"0.443735"
$ws.Cells.Item(12, 8).Value2 = [double]"0.88747"
$ws.Cells.Item(12, 9).Value2 = [double]"0.02857803451828042"
$ws.Cells.Item(12, 10).Value2 = [double]"0.01923525830411462"
$ws.Cells.Item(12, 11).Value2 = [double]"2"
$ws.Cells.Item(12, 12).Value2 = [double]"1"
$ws.Cells.Item(12, 13).Value2 = [double]"11.9507005"
$ws.Cells.Item(12, 14).Value2 = [double]"23.901401"
$ws.Cells.Item(12, 15).Value2 = [double]"0.2543137660693869"
$ws.Cells.Item(12, 16).Value2 = [double]"0.1852460510065796"
$ws.Cells.Item(12, 17).Value2 = [double]"5.3029440863675"
$ws.Cells.Item(12, 18).Value2 = [double]"21.21177634547"
$ws.Cells.Item(12, 19).Value2 = [double]"0.007267787585204831"
$ws.Cells.Item(12, 20).Value2 = [double]"0.003563255640928752"

# Row 13
$ws.Cells.Item(13, 1).Value2 = "MuSCs"
$ws.Cells.Item(13, 2).Value2 = "Myoc"
$ws.Cells.Item(13, 3).Value2 = "Fzd4"
$ws.Cells.Item(13, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value2 = [double]"1"
$ws.Cells.Item(13, 6).Value2 = [double]"0.5"
$ws.Cells.Item(13, 7).Value2 = [double]"0.443735"
$ws.Cells.Item(13, 8).Value2 = [double]"0.88747"
$ws.Cells.Item(13, 9).Value2 = [double]"0.02857803451828042"
$ws.Cells.Item(13, 10).Value2 = [double]"0.01923525830411462"
$ws.Cells.Item(13, 11).Value2 = [double]"2"
$ws.Cells.Item(13, 12).Value2 = [double]"0.6666666666666666"
$ws.Cells.Item(13, 13).Value2 = [double]"0.272275"
$ws.Cells.Item(13, 14).Value2 = [double]"0.816825"
$ws.Cells.Item(13, 15).Value2 = [double]"0.005794077146903843"
$ws.Cells.Item(13, 16).Value2 = [double]"0.006330742102249548"
$ws.Cells.Item(13, 17).Value2 = [double]"0.120817947125"
$ws.Cells.Item(13, 18).Value2 = [double]"0.7249076827500001"
$ws.Cells.Item(13, 19).Value2 = [double]"0.0001655833367057978"
$ws.Cells.Item(13, 20).Value2 = [double]"0.0001217734595935037"

